$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "21.663.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.534.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9999"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.60"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3946"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3168"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07161"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.061"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.08%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.705"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.34"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.609"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.539.16"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001093"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06601"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.96"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.140"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.52"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.67"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.362"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "21.639.88"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.351"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.98"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.35"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.849"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.709.04"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.27"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.072"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9430"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -11.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08146"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.177"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.482"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -8.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05986"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02214"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.456"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -11.35%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.13"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.27%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2034"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.00%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.181"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.69%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9995"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5787"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.01"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.719"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5548"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.168"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.878"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "116.09"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06688"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.44%  "
